# Auto-generated Excel COM-interop script
# Applies the Gilgamesh_Profits.xlsx commit: updates currentAveragePrice /
# LevePrice / LeveProfit columns (H-N) for a set of rows across several
# worksheets (ALC, ARM, CRP, CUL, GSM, LTW, WVR), as produced by the
# scheduled market-data refresh runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(42, 8).Value = 317.33334  # H42
$ws.Cells.Item(42, 9).Value = 133.5  # I42
$ws.Cells.Item(42, 10).Value = 409.25  # J42
$ws.Cells.Item(42, 11).Value = 400.5  # K42
$ws.Cells.Item(42, 12).Value = 1227.75  # L42
$ws.Cells.Item(42, 13).Value = -170.5  # M42
$ws.Cells.Item(42, 14).Value = -1687.75  # N42

$ws.Cells.Item(58, 8).Value = 249.33333  # H58
$ws.Cells.Item(58, 9).Value = 109.2  # I58
$ws.Cells.Item(58, 10).Value = 950  # J58
$ws.Cells.Item(58, 11).Value = 327.6  # K58
$ws.Cells.Item(58, 12).Value = 2850  # L58
$ws.Cells.Item(58, 13).Value = -177.6  # M58
$ws.Cells.Item(58, 14).Value = -3150  # N58

$ws.Cells.Item(61, 8).Value = 292.6  # H61
$ws.Cells.Item(61, 9).Value = 292.6  # I61
$ws.Cells.Item(61, 11).Value = 877.8000000000001  # K61
$ws.Cells.Item(61, 13).Value = -705.8000000000001  # M61

$ws.Cells.Item(82, 8).Value = 677.3333  # H82
$ws.Cells.Item(82, 9).Value = 677.3333  # I82
$ws.Cells.Item(82, 11).Value = 2031.9999  # K82
$ws.Cells.Item(82, 13).Value = -1625.9999  # M82

$ws.Cells.Item(85, 8).Value = 677.3333  # H85
$ws.Cells.Item(85, 9).Value = 677.3333  # I85
$ws.Cells.Item(85, 11).Value = 2031.9999  # K85
$ws.Cells.Item(85, 13).Value = -627.9999  # M85

$ws.Cells.Item(99, 8).Value = 454.08334  # H99
$ws.Cells.Item(99, 10).Value = 1125  # J99
$ws.Cells.Item(99, 12).Value = 3375  # L99
$ws.Cells.Item(99, 14).Value = -6371  # N99

$ws.Cells.Item(104, 8).Value = 343.5  # H104
$ws.Cells.Item(104, 9).Value = 212.4  # I104
$ws.Cells.Item(104, 11).Value = 637.2  # K104
$ws.Cells.Item(104, 13).Value = 1109.8  # M104

$ws.Cells.Item(118, 8).Value = 693.25  # H118
$ws.Cells.Item(118, 9).Value = 652.4  # I118
$ws.Cells.Item(118, 10).Value = 897.5  # J118
$ws.Cells.Item(118, 11).Value = 1957.2  # K118
$ws.Cells.Item(118, 12).Value = 2692.5  # L118
$ws.Cells.Item(118, 13).Value = -300.1999999999998  # M118
$ws.Cells.Item(118, 14).Value = -6006.5  # N118

$ws.Cells.Item(127, 8).Value = 566.5  # H127
$ws.Cells.Item(127, 9).Value = 566.5  # I127
$ws.Cells.Item(127, 11).Value = 1699.5  # K127
$ws.Cells.Item(127, 13).Value = 3260.5  # M127

$ws.Cells.Item(129, 8).Value = 2410.2666  # H129
$ws.Cells.Item(129, 9).Value = 1063.1666  # I129
$ws.Cells.Item(129, 10).Value = 3308.3333  # J129
$ws.Cells.Item(129, 11).Value = 3189.4998  # K129
$ws.Cells.Item(129, 12).Value = 9924.999899999999  # L129
$ws.Cells.Item(129, 13).Value = 1810.5002  # M129
$ws.Cells.Item(129, 14).Value = -19924.9999  # N129

$ws.Cells.Item(138, 8).Value = 3588.5  # H138
$ws.Cells.Item(138, 9).Value = 3181.348  # I138
$ws.Cells.Item(138, 11).Value = 9544.044  # K138
$ws.Cells.Item(138, 13).Value = -4404.044  # M138

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4268.875  # H32
$ws.Cells.Item(32, 9).Value = 4237.1304  # I32
$ws.Cells.Item(32, 11).Value = 4237.1304  # K32
$ws.Cells.Item(32, 13).Value = -3950.1304  # M32

$ws.Cells.Item(74, 8).Value = 200300.42  # H74
$ws.Cells.Item(74, 9).Value = 371596  # I74
$ws.Cells.Item(74, 10).Value = 2651.6924  # J74
$ws.Cells.Item(74, 11).Value = 371596  # K74
$ws.Cells.Item(74, 12).Value = 2651.6924  # L74
$ws.Cells.Item(74, 13).Value = -370722  # M74
$ws.Cells.Item(74, 14).Value = -4399.6924  # N74

$ws.Cells.Item(77, 8).Value = 200300.42  # H77
$ws.Cells.Item(77, 9).Value = 371596  # I77
$ws.Cells.Item(77, 10).Value = 2651.6924  # J77
$ws.Cells.Item(77, 11).Value = 1857980  # K77
$ws.Cells.Item(77, 12).Value = 13258.462  # L77
$ws.Cells.Item(77, 13).Value = -1853612  # M77
$ws.Cells.Item(77, 14).Value = -21994.462  # N77

$ws.Cells.Item(132, 8).Value = 3064.5227  # H132
$ws.Cells.Item(132, 9).Value = 2828.861  # I132
$ws.Cells.Item(132, 11).Value = 8486.582999999999  # K132
$ws.Cells.Item(132, 13).Value = -5956.582999999999  # M132

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1667.4706  # H16
$ws.Cells.Item(16, 9).Value = 1708.3077  # I16
$ws.Cells.Item(16, 10).Value = 1534.75  # J16
$ws.Cells.Item(16, 11).Value = 1708.3077  # K16
$ws.Cells.Item(16, 12).Value = 1534.75  # L16
$ws.Cells.Item(16, 13).Value = -1421.3077  # M16
$ws.Cells.Item(16, 14).Value = -2108.75  # N16

$ws.Cells.Item(31, 8).Value = 3524.4792  # H31
$ws.Cells.Item(31, 9).Value = 2645.1282  # I31
$ws.Cells.Item(31, 11).Value = 2645.1282  # K31
$ws.Cells.Item(31, 13).Value = -2350.1282  # M31

$ws.Cells.Item(34, 8).Value = 3524.4792  # H34
$ws.Cells.Item(34, 9).Value = 2645.1282  # I34
$ws.Cells.Item(34, 11).Value = 2645.1282  # K34
$ws.Cells.Item(34, 13).Value = -2443.1282  # M34

$ws.Cells.Item(58, 8).Value = 2640.6667  # H58
$ws.Cells.Item(58, 9).Value = 1387.1428  # I58
$ws.Cells.Item(58, 10).Value = 3737.5  # J58
$ws.Cells.Item(58, 11).Value = 1387.1428  # K58
$ws.Cells.Item(58, 12).Value = 3737.5  # L58
$ws.Cells.Item(58, 13).Value = -1184.1428  # M58
$ws.Cells.Item(58, 14).Value = -4143.5  # N58

$ws.Cells.Item(113, 8).Value = 1667.4706  # H113
$ws.Cells.Item(113, 9).Value = 1708.3077  # I113
$ws.Cells.Item(113, 10).Value = 1534.75  # J113
$ws.Cells.Item(113, 11).Value = 1708.3077  # K113
$ws.Cells.Item(113, 12).Value = 1534.75  # L113
$ws.Cells.Item(113, 13).Value = 461.6922999999999  # M113
$ws.Cells.Item(113, 14).Value = -5874.75  # N113

$ws.Cells.Item(120, 8).Value = 36331  # H120
$ws.Cells.Item(120, 10).Value = 36331  # J120
$ws.Cells.Item(120, 12).Value = 36331  # L120
$ws.Cells.Item(120, 14).Value = -43589  # N120

$ws.Cells.Item(121, 8).Value = 49326  # H121
$ws.Cells.Item(121, 10).Value = 49326  # J121
$ws.Cells.Item(121, 12).Value = 49326  # L121
$ws.Cells.Item(121, 14).Value = -51946  # N121

$ws.Cells.Item(134, 8).Value = 2993.8948  # H134
$ws.Cells.Item(134, 9).Value = 2758.9333  # I134
$ws.Cells.Item(134, 11).Value = 8276.7999  # K134
$ws.Cells.Item(134, 13).Value = -5741.7999  # M134

$ws.Cells.Item(136, 8).Value = 2640.6667  # H136
$ws.Cells.Item(136, 9).Value = 1387.1428  # I136
$ws.Cells.Item(136, 10).Value = 3737.5  # J136
$ws.Cells.Item(136, 11).Value = 4161.428400000001  # K136
$ws.Cells.Item(136, 12).Value = 11212.5  # L136
$ws.Cells.Item(136, 13).Value = -1611.428400000001  # M136
$ws.Cells.Item(136, 14).Value = -16312.5  # N136

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 63934840  # H4
$ws.Cells.Item(4, 9).Value = 74617180  # I4
$ws.Cells.Item(4, 11).Value = 223851540  # K4
$ws.Cells.Item(4, 13).Value = -223851428  # M4

$ws.Cells.Item(6, 8).Value = 119.8  # H6
$ws.Cells.Item(6, 9).Value = 124.75  # I6
$ws.Cells.Item(6, 10).Value = 100  # J6
$ws.Cells.Item(6, 11).Value = 374.25  # K6
$ws.Cells.Item(6, 12).Value = 300  # L6
$ws.Cells.Item(6, 13).Value = -261.25  # M6
$ws.Cells.Item(6, 14).Value = -526  # N6

$ws.Cells.Item(10, 8).Value = 27778044  # H10
$ws.Cells.Item(10, 9).Value = 27778044  # I10
$ws.Cells.Item(10, 11).Value = 83334132  # K10
$ws.Cells.Item(10, 13).Value = -83333993  # M10

$ws.Cells.Item(11, 8).Value = 854.0454999999999  # H11
$ws.Cells.Item(11, 10).Value = 700  # J11
$ws.Cells.Item(11, 12).Value = 2100  # L11
$ws.Cells.Item(11, 14).Value = -2380  # N11

$ws.Cells.Item(15, 8).Value = 2366.6667  # H15
$ws.Cells.Item(15, 9).Value = 100  # I15
$ws.Cells.Item(15, 11).Value = 300  # K15
$ws.Cells.Item(15, 13).Value = -160  # M15

$ws.Cells.Item(17, 8).Value = 27777958  # H17
$ws.Cells.Item(17, 10).Value = 170  # J17
$ws.Cells.Item(17, 12).Value = 510  # L17
$ws.Cells.Item(17, 14).Value = -848  # N17

$ws.Cells.Item(24, 8).Value = 166668850  # H24
$ws.Cells.Item(24, 10).Value = 2962.25  # J24
$ws.Cells.Item(24, 12).Value = 8886.75  # L24
$ws.Cells.Item(24, 14).Value = -9346.75  # N24

$ws.Cells.Item(39, 8).Value = 6466.909  # H39
$ws.Cells.Item(39, 10).Value = 6924.8  # J39
$ws.Cells.Item(39, 12).Value = 20774.4  # L39
$ws.Cells.Item(39, 14).Value = -21362.4  # N39

$ws.Cells.Item(47, 8).Value = 1904.6  # H47
$ws.Cells.Item(47, 9).Value = 438.42856  # I47
$ws.Cells.Item(47, 11).Value = 1315.28568  # K47
$ws.Cells.Item(47, 13).Value = -884.28568  # M47

$ws.Cells.Item(51, 8).Value = 425  # H51
$ws.Cells.Item(51, 10).Value = 350  # J51
$ws.Cells.Item(51, 12).Value = 1050  # L51
$ws.Cells.Item(51, 14).Value = -1970  # N51

$ws.Cells.Item(57, 14).ClearContents()  # N57
$ws.Cells.Item(57, 8).Value = 1998.5  # H57
$ws.Cells.Item(57, 9).Value = 1998.5  # I57
$ws.Cells.Item(57, 10).Value = 0  # J57
$ws.Cells.Item(57, 11).Value = 5995.5  # K57
$ws.Cells.Item(57, 12).Value = 0  # L57
$ws.Cells.Item(57, 13).Value = -5436.5  # M57

$ws.Cells.Item(59, 8).Value = 334.5  # H59
$ws.Cells.Item(59, 9).Value = 334.5  # I59
$ws.Cells.Item(59, 11).Value = 1003.5  # K59
$ws.Cells.Item(59, 13).Value = -463.5  # M59

$ws.Cells.Item(101, 8).Value = 3666.3333  # H101
$ws.Cells.Item(101, 9).Value = 1499.5  # I101
$ws.Cells.Item(101, 10).Value = 8000  # J101
$ws.Cells.Item(101, 11).Value = 4498.5  # K101
$ws.Cells.Item(101, 12).Value = 24000  # L101
$ws.Cells.Item(101, 13).Value = -2064.5  # M101
$ws.Cells.Item(101, 14).Value = -28868  # N101

$ws.Cells.Item(107, 8).Value = 279.91666  # H107
$ws.Cells.Item(107, 10).Value = 279.91666  # J107
$ws.Cells.Item(107, 12).Value = 839.7499799999999  # L107
$ws.Cells.Item(107, 14).Value = -4679.74998  # N107

$ws.Cells.Item(116, 8).Value = 132091.28  # H116
$ws.Cells.Item(116, 9).Value = 301833  # I116
$ws.Cells.Item(116, 10).Value = 4785  # J116
$ws.Cells.Item(116, 11).Value = 905499  # K116
$ws.Cells.Item(116, 12).Value = 14355  # L116
$ws.Cells.Item(116, 13).Value = -902057  # M116
$ws.Cells.Item(116, 14).Value = -21239  # N116

$ws.Cells.Item(129, 8).Value = 1830.1428  # H129
$ws.Cells.Item(129, 10).Value = 1627.6666  # J129
$ws.Cells.Item(129, 12).Value = 4882.9998  # L129
$ws.Cells.Item(129, 14).Value = -14882.9998  # N129

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 2566.6  # H132
$ws.Cells.Item(132, 9).Value = 2049.9167  # I132
$ws.Cells.Item(132, 11).Value = 6149.750100000001  # K132
$ws.Cells.Item(132, 13).Value = -3619.750100000001  # M132

$ws.Cells.Item(135, 8).Value = 68691.53999999999  # H135
$ws.Cells.Item(135, 10).Value = 68691.53999999999  # J135
$ws.Cells.Item(135, 12).Value = 68691.53999999999  # L135
$ws.Cells.Item(135, 14).Value = -78831.53999999999  # N135

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1147.8462  # H22
$ws.Cells.Item(22, 10).Value = 656.8333  # J22
$ws.Cells.Item(22, 12).Value = 656.8333  # L22
$ws.Cells.Item(22, 14).Value = -1246.8333  # N22

$ws.Cells.Item(27, 8).Value = 1147.8462  # H27
$ws.Cells.Item(27, 10).Value = 656.8333  # J27
$ws.Cells.Item(27, 12).Value = 656.8333  # L27
$ws.Cells.Item(27, 14).Value = -870.8333  # N27

$ws.Cells.Item(68, 8).Value = 2817  # H68
$ws.Cells.Item(68, 9).Value = 3399  # I68
$ws.Cells.Item(68, 10).Value = 2235  # J68
$ws.Cells.Item(68, 11).Value = 3399  # K68
$ws.Cells.Item(68, 12).Value = 2235  # L68
$ws.Cells.Item(68, 13).Value = -2650  # M68
$ws.Cells.Item(68, 14).Value = -3733  # N68

$ws.Cells.Item(71, 8).Value = 2817  # H71
$ws.Cells.Item(71, 9).Value = 3399  # I71
$ws.Cells.Item(71, 10).Value = 2235  # J71
$ws.Cells.Item(71, 11).Value = 16995  # K71
$ws.Cells.Item(71, 12).Value = 11175  # L71
$ws.Cells.Item(71, 13).Value = -13251  # M71
$ws.Cells.Item(71, 14).Value = -18663  # N71

$ws.Cells.Item(132, 8).Value = 4203.075  # H132
$ws.Cells.Item(132, 10).Value = 8035.2856  # J132
$ws.Cells.Item(132, 12).Value = 24105.8568  # L132
$ws.Cells.Item(132, 14).Value = -29165.8568  # N132

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(74, 8).Value = 19750  # H74
$ws.Cells.Item(74, 10).Value = 19750  # J74
$ws.Cells.Item(74, 12).Value = 19750  # L74
$ws.Cells.Item(74, 14).Value = -21622  # N74

$ws.Cells.Item(77, 8).Value = 19750  # H77
$ws.Cells.Item(77, 10).Value = 19750  # J77
$ws.Cells.Item(77, 12).Value = 59250  # L77
$ws.Cells.Item(77, 14).Value = -68610  # N77

$ws.Cells.Item(136, 8).Value = 23629.75  # H136
$ws.Cells.Item(136, 9).Value = 22099.375  # I136
$ws.Cells.Item(136, 11).Value = 66298.125  # K136
$ws.Cells.Item(136, 13).Value = -63748.125  # M136
